$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "60.333.14"
Set-TextValue "E2" "  -1.05%  "
Set-TextValue "D3" "2.612.68"
Set-TextValue "E3" "  +0.16%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "585.91"
Set-TextValue "E5" "  +2.55%  "
Set-TextValue "D6" "143.24"
Set-TextValue "E6" "  -0.10%  "
Set-TextValue "E7" "  +0.21%  "
Set-TextValue "D8" "0.597"
Set-TextValue "E8" "  -0.77%  "
Set-TextValue "D9" "6.50"
Set-TextValue "E9" "  -0.10%  "
Set-TextValue "E10" "  -1.33%  "
Set-TextValue "D11" "0.155"
Set-TextValue "E11" "  +1.02%  "
Set-TextValue "E12" "  +1.32%  "
Set-TextValue "D13" "3.074.22"
Set-TextValue "E13" "  -0.02%  "
Set-TextValue "D14" "24.81"
Set-TextValue "E14" "  +5.39%  "
Set-TextValue "D15" "60.328.52"
Set-TextValue "E15" "  -1.04%  "
Set-TextValue "E16" "  -0.58%  "
Set-TextValue "D17" "2.619.52"
Set-TextValue "E17" "  -0.12%  "
Set-TextValue "D18" "11.41"
Set-TextValue "E19" "  -0.37%  "
Set-TextValue "D20" "346.39"
Set-TextValue "E20" "  -1.09%  "
Set-TextValue "D21" "6.91"
Set-TextValue "E21" "  -2.77%  "
Set-TextValue "E22" "  -0.21%  "
Set-TextValue "D23" "0.536"
Set-TextValue "E23" "  +3.04%  "
Set-TextValue "D24" "63.77"
Set-TextValue "E24" "  -0.81%  "
Set-TextValue "E25" "  +0.41%  "
Set-TextValue "E26" "  -0.29%  "
Set-TextValue "D27" "8.02"
Set-TextValue "E27" "  +4.01%  "
Set-TextValue "D28" "1.93"
Set-TextValue "E28" "  +5.31%  "
Set-TextValue "D29" "0.0₃0797"
Set-TextValue "E29" "  -0.01%  "
Set-TextValue "E30" "  +1.73%  "
Set-TextValue "D31" "168.53"
Set-TextValue "E31" "  +5.18%  "
Set-TextValue "D32" "0.999"
Set-TextValue "E32" "  +0.17%  "
Set-TextValue "E33" "  -0.30%  "
Set-TextValue "E34" "  +5.41%  "
Set-TextValue "E35" "  +0.58%  "
Set-TextValue "E36" "  +7.90%  "
Set-TextValue "E37" "  +2.31%  "
Set-TextValue "D38" "319.86"
Set-TextValue "E38" "  +6.92%  "
Set-TextValue "D39" "38.40"
Set-TextValue "E39" "  +1.56%  "
Set-TextValue "E40" "  +3.19%  "
Set-TextValue "D41" "0.848"
Set-TextValue "E41" "  -0.72%  "
Set-TextValue "D42" "135.62"
Set-TextValue "E42" "  -3.28%  "
Set-TextValue "E43" "  +0.28%  "
Set-TextValue "D44" "0.999"
Set-TextValue "E44" "  +0.30%  "
Set-TextValue "D45" "19.95"
Set-TextValue "E45" "  +2.05%  "
Set-TextValue "E46" "  +0.07%  "
Set-TextValue "E47" "  +1.77%  "
Set-TextValue "E48" "  +0.28%  "
Set-TextValue "D49" "20.07"
Set-TextValue "E49" "  +1.90%  "
Set-TextValue "D50" "0.0241"
Set-TextValue "E50" "  +0.01%  "
Set-TextValue "D51" "10.73"
Set-TextValue "E51" "  +0.33%  "
